# Updates the "cryptos" price/volume table (columns D and E, rows 2-51)
# to the refreshed values from the latest GitHub Actions run.
#
# Some Price-column values (e.g. "0.998", "24.95") are strings that also
# happen to look like plain numbers. Assigning them straight to .Value
# would let Excel auto-convert the cell to a real number (dropping
# significant trailing zeros such as "12.50" -> 12.5). To keep them as
# literal text - matching the source data - they're written with a
# leading apostrophe (forces text entry, like typing into the grid) and
# then the resulting "quote prefix" cell format is cleared back to the
# default Normal style so only the cell's value/type changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.438.79"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "3.137.38"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'572.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'164.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.574"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.92%  "
$ws.Range("D9").Value = "3.151.97"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("E11").Value = "  -3.30%  "
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "3.686.42"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").Value = "'0.128"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.49%  "
$ws.Range("D15").Value = "64.454.98"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'24.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").Value = "3.150.54"
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("E18").Value = "  -2.38%  "
$ws.Range("D19").Value = "'407.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("D20").Value = "'12.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.93%  "
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("D22").Value = "'7.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'68.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").Value = "'0.483"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.63%  "
$ws.Range("E26").Value = "  -5.49%  "
$ws.Range("E27").Value = "  -3.18%  "
$ws.Range("D28").Value = "'8.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "'1.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").Value = "'21.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.79%  "
$ws.Range("D33").Value = "'163.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.78%  "
$ws.Range("D34").Value = "'4.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("D35").Value = "'6.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").Value = "'1.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("D39").Value = "2.632.22"
$ws.Range("E39").Value = "  -2.91%  "
$ws.Range("D40").Value = "'23.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").Value = "'4.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("D42").Value = "'38.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("D43").Value = "'0.690"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.70%  "
$ws.Range("D44").Value = "'0.0612"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("D45").Value = "'5.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.04%  "
$ws.Range("D46").Value = "'289.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "'21.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Value = "'0.0253"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.53%  "
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").Value = "'0.0972"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("D51").Value = "'10.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.18%  "
